$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.637.14"
$ws.Range("E2").Value = "  +3.58%  "
$ws.Range("D3").Value = "1.851.88"
$ws.Range("E3").Value = "  +2.69%  "
$ws.Range("D4").Value = "'1.029"
$ws.Range("E4").Value = "  +2.99%  "
$ws.Range("D5").Value = "'320.25"
$ws.Range("E5").Value = "  +3.68%  "
$ws.Range("D6").Value = "'1.028"
$ws.Range("E6").Value = "  +3.03%  "
$ws.Range("D7").Value = "'0.4380"
$ws.Range("E7").Value = "  +2.20%  "
$ws.Range("D8").Value = "'0.3751"
$ws.Range("E8").Value = "  +3.07%  "
$ws.Range("D9").Value = "'0.07417"
$ws.Range("E9").Value = "  +3.14%  "
$ws.Range("D10").Value = "'0.8773"
$ws.Range("E10").Value = "  +1.90%  "
$ws.Range("D11").Value = "'21.54"
$ws.Range("E11").Value = "  +3.90%  "
$ws.Range("D12").Value = "1.851.61"
$ws.Range("E12").Value = "  -1.97%  "
$ws.Range("D13").Value = "'5.495"
$ws.Range("E13").Value = "  +3.53%  "
$ws.Range("D14").Value = "'6.695"
$ws.Range("E14").Value = "  +1.77%  "
$ws.Range("D15").Value = "'0.07161"
$ws.Range("E15").Value = "  +3.79%  "
$ws.Range("D16").Value = "'82.94"
$ws.Range("E16").Value = "  +3.81%  "
$ws.Range("D17").Value = "'1.033"
$ws.Range("E17").Value = "  +3.43%  "
$ws.Range("D18").Value = "'0.000009033"
$ws.Range("E18").Value = "  +2.03%  "
$ws.Range("D19").Value = "'1.027"
$ws.Range("E19").Value = "  +2.56%  "
$ws.Range("D20").Value = "'15.45"
$ws.Range("E20").Value = "  +1.89%  "
$ws.Range("D21").Value = "27.616.67"
$ws.Range("E21").Value = "  +3.53%  "
$ws.Range("D22").Value = "'5.256"
$ws.Range("E22").Value = "  +1.97%  "
$ws.Range("E23").Value = "  +1.23%  "
$ws.Range("D24").Value = "2.066.69"
$ws.Range("E24").Value = "  -2.03%  "
$ws.Range("D25").Value = "'157.53"
$ws.Range("E25").Value = "  +3.62%  "
$ws.Range("D26").Value = "'1.930"
$ws.Range("E26").Value = "  +5.01%  "
$ws.Range("D27").Value = "'18.75"
$ws.Range("E27").Value = "  +2.99%  "
$ws.Range("D28").Value = "'5.270"
$ws.Range("E28").Value = "  +1.74%  "
$ws.Range("D29").Value = "'1.947"
$ws.Range("E29").Value = "  +2.88%  "
$ws.Range("D30").Value = "'116.18"
$ws.Range("E30").Value = "  +1.05%  "
$ws.Range("D31").Value = "'0.09086"
$ws.Range("E31").Value = "  +1.85%  "
$ws.Range("D32").Value = "'1.208"
$ws.Range("E32").Value = "  +4.12%  "
$ws.Range("D33").Value = "'0.7684"
$ws.Range("E33").Value = "  +2.61%  "
$ws.Range("D34").Value = "'4.521"
$ws.Range("E34").Value = "  +2.92%  "
$ws.Range("D35").Value = "'2.881"
$ws.Range("E35").Value = "  +4.70%  "
$ws.Range("D36").Value = "'1.030"
$ws.Range("E36").Value = "  +2.93%  "
$ws.Range("E37").Value = "  +2.24%  "
$ws.Range("D38").Value = "'0.01981"
$ws.Range("E38").Value = "  +4.02%  "
$ws.Range("D39").Value = "'0.05278"
$ws.Range("E39").Value = "  +2.06%  "
$ws.Range("D40").Value = "'0.5168"
$ws.Range("E40").Value = "  +3.35%  "
$ws.Range("D41").Value = "'2.802"
$ws.Range("E41").Value = "  +7.02%  "
$ws.Range("D42").Value = "'0.1676"
$ws.Range("E42").Value = "  +2.95%  "
$ws.Range("D43").Value = "'6.693"
$ws.Range("E43").Value = "  +3.78%  "
$ws.Range("D44").Value = "'8.570"
$ws.Range("E44").Value = "  +4.09%  "
$ws.Range("D45").Value = "'108.99"
$ws.Range("E45").Value = "  +2.65%  "
$ws.Range("D46").Value = "'10.57"
$ws.Range("E46").Value = "  +3.17%  "
$ws.Range("D47").Value = "'1.717"
$ws.Range("E47").Value = "  +4.33%  "
$ws.Range("D48").Value = "'0.4655"
$ws.Range("E48").Value = "  +2.41%  "
$ws.Range("D49").Value = "'0.06381"
$ws.Range("E49").Value = "  +2.34%  "
$ws.Range("D50").Value = "'1.890"
$ws.Range("E50").Value = "  +5.67%  "
$ws.Range("D51").Value = "'39.63"
$ws.Range("E51").Value = "  +6.79%  "